$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$inc  = $wb.Worksheets.Item("Include from SNOMED CT")

# Replicate the formatting (border/fill/alignment) of the last existing
# metadata row down into the new rows (11-16) before writing any values,
# so the freshly-created rows pick up style index 2 like their neighbours.
$meta.Range("A10:B10").Copy()
$meta.Range("A11:B16").PasteSpecial(-4122)

# --- Metadata sheet value updates ---

# Version: 0.1.6 -> 0.1.7
$meta.Range("B3").Value = "0.1.7"

# Status: active -> draft
$meta.Range("B6").Value = "draft"

# Date
$meta.Range("B8").Value = "2024-11-22T12:33:30-06:00"

# Contact (row 10): "No display for ContactDetail" -> publisher org + url
$meta.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"

# New row 11: second Contact entry
$meta.Range("A11").Value = "Contact"
$meta.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# New row 12: Jurisdiction (blank value)
$meta.Range("A12").Value = "Jurisdiction"
$meta.Range("B12").Value = ""

# Row 13: Description
$meta.Range("A13").Value = "Description"
$meta.Range("B13").Value = "The myelodysplastic syndromes (MDS) are a group of clonal hematopoietic stem cell diseases characterized by cytopenia(s), dysplasia (abnormal growth or development leading to an alteration in size, shape, and organization of the cell) in one or more of the major myeloid cell lines (WBC, RBC, and/or platelets), ineffective hematopoiesis, and an increased risk of developing acute myelogenous leukemia (AML). MDS occurs primarily in older adults, with a median age of 70 years. The majority of recipients present with symptoms related to cytopenias. Most recipients present with anemia requiring RBC transfusions."

# Row 14: Purpose (blank value)
$meta.Range("A14").Value = "Purpose"
$meta.Range("B14").Value = ""

# Row 15: Copyright (blank value)
$meta.Range("A15").Value = "Copyright"
$meta.Range("B15").Value = ""

# New row 16: Immutable
$meta.Range("A16").Value = "Immutable"
$meta.Range("B16").Value = "BooleanType[null]"

# --- Include from SNOMED CT sheet updates ---
# Row 3 (previously the blank-display concept row) becomes fully blank.
$inc.Range("A3").Value = ""
$inc.Range("B3").Value = ""
